# Apply the "Add files via upload" edit:
#   - Rename Sheet1 -> Report
#   - Update membership counts (column B) for each society
#   - Re-enter the IOS name in A4 (forces a fresh shared-string entry,
#     matching the reshuffled sharedStrings.xml in the target workbook)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "Report"

# Update the membership-count column (B) for each society row.
$ws.Range("B2").Value = 600   # FLASCO
$ws.Range("B3").Value = 500   # GASCO
$ws.Range("B4").Value = 200   # IOS (Indiana Oncology Society)
$ws.Range("B5").Value = 176   # IOWA Oncology Society
$ws.Range("B6").Value = 400   # MOASC

# Re-type the IOS society name so it is re-interned as a shared string,
# mirroring the source workbook's change.
$ws.Range("A4").Value = "IOS (Indiana Oncology Society)"

# Move the active selection, as recorded in the saved view state.
$ws.Range("D12").Select()
